$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-Cell($row, $col, $new) {
    # Assign the cell Range's Text directly. Table cells hold a single
    # run/value in this report, so this swaps the number while leaving the
    # run/paragraph formatting (rPr/pPr/tcPr) untouched - and, unlike
    # Find.Execute with Replace:=wdReplaceAll, it only ever touches this
    # one cell (important since short values like "5.5"/"2.5" also occur
    # as substrings of other, unrelated cell values elsewhere in the
    # table, e.g. "12 (25.5)").
    $t.Cell($row, $col).Range.Text = $new
}

Set-Cell 2  5 "12 ( 9.4)"
Set-Cell 2  7 "7 (5.5)"
Set-Cell 3  5 "28 (22.0)"
Set-Cell 4  5 "41 (32.3)"
Set-Cell 5  3 "5.0"
Set-Cell 5  5 "52 (40.9)"
Set-Cell 5  7 "6 (4.7)"
Set-Cell 6  5 "56 (44.1)"
Set-Cell 7  3 "2.0"
Set-Cell 7  5 "62 (48.8)"
Set-Cell 8  5 "69 (54.3)"
Set-Cell 9  5 "67 (52.8)"
Set-Cell 10 5 "72 (56.7)"
Set-Cell 11 5 "75 (59.1)"
Set-Cell 12 5 "77 (60.6)"
Set-Cell 13 5 "82 (64.6)"
